# Commit: "Finito: questa e la versione il cui codice finisce sul sito,
# tutta commentata bene e con pieno utilizzo di contact (ovviamente non va!)"
#
# The "Versione senza contact" section (Sensore / UserCmd alternative
# writeup at the end of the doc) gets re-highlighted from lightGray to
# cyan, its two sub-headings ("Scontrol" / "Riferimenti e configurazione")
# pick up the same cyan highlight (including the paragraph-mark run
# properties), and the trailing ". Lasci comunque ... avvio" + "."
# runs get merged into a single highlighted run.

$d = $word.ActiveDocument

# First, merge the trailing ". Lasci comunque invariato in Edi
# l'assegnazione di nomi ai processi e il loro avvio" run with the lone
# "." run that follows it, by replacing the combined text with itself —
# Word collapses the matched span into one run using the first run's
# formatting.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Text = "il loro avvio."
$find.Replacement.Text = "il loro avvio."
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)

$paras = $d.Paragraphs

# "Sensore" sub-section body (Nel metodo sendData ... contact)
$paras.Item(40).Range.HighlightColorIndex = "Turquoise"

# "UserCmd" sub-section body (nel metodo MandaComando ... restituito)
$paras.Item(42).Range.HighlightColorIndex = "Turquoise"

# "Scontrol" sub-heading — also colors the paragraph-mark run properties
$paras.Item(43).Range.Font.HighlightColorIndex = "Turquoise"

# Body paragraph under "Scontrol" (Puoi commentare ... normali chiamate)
# this paragraph previously had no <w:pPr> at all; setting via Font
# creates one with the highlighted paragraph mark, matching the diff.
$paras.Item(44).Range.Font.HighlightColorIndex = "Turquoise"

# "Riferimenti e configurazione" sub-heading
$paras.Item(45).Range.Font.HighlightColorIndex = "Turquoise"

# Body paragraph (Sia i sensori che UserCmd ... il loro avvio.)
$paras.Item(46).Range.HighlightColorIndex = "Turquoise"
